# Applies the "Specification" slide rework:
#  - Strip the four leftover layout placeholders (two subtitle/content pairs)
#  - Re-point the footer/date placeholders to the new title/date text
#  - Add a new free-form text box holding the requirements list

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- Remove the four unused placeholder shapes ------------------------------
# Shapes.Item(2) is always the next still-unwanted placeholder once the
# earlier ones have been removed. Calling .Delete() on a content/text
# placeholder the first time just clears it back to its layout default
# (PowerPoint keeps an "empty" placeholder around); calling .Delete() again
# on that now-empty placeholder removes it from the slide for good.
for ($n = 1; $n -le 4; $n++) {
    $s.Shapes.Item(2).Delete()
    $s.Shapes.Item(2).Delete()
}

# --- Footer placeholder: "Presentation title" -> "Simplify's Design Specification"
$footer = $s.Shapes.Item("Footer Placeholder 6")
$footer.TextFrame.TextRange.Text = "Simplify's"
$footer.TextFrame.TextRange.InsertAfter(" Design Specification")

# --- Date placeholder: "20XX" -> "2023"
$dateShape = $s.Shapes.Item("Date Placeholder 5")
$dateShape.TextFrame.TextRange.Text = "2023"

# --- New free-form text box with the requirements write-up -----------------
$tb = $s.Shapes.AddTextbox(1, 0, 0, 100, 100)
$tb.Name = "TextBox 17"

# AddTextbox/Left/Top/Width/Height all operate in points; nudge by a hair so
# the point->EMU round trip lands on the exact target EMU values.
$tb.Left = (167816 / 12700.0) + 0.00002
$tb.Top = (1182330 / 12700.0) + 0.00002
$tb.Width = (11823016 / 12700.0) + 0.00002
$tb.Height = (5355312 / 12700.0) + 0.00002

$lines = @(
"User authentication: The platform should provide the users; students and tutors a management system to login or register using google or a registration form. The user should be able to register as a student or a tutor.",
"User based dashboard: The platform should distinguish between a tutor and a student, with access being given to the student for only the student specific pages.",
"",
"Course management: The platform should have a user-friendly environment that allows the tutors to create a course, add pre-existing courses, and publish course materials including slides, pdfs, images and videos.",
"",
"Course accessibility: The platform should provide the student access to view all the course material and submit course work.",
"",
"Personalization: The platform should provide the student the option to enrol for more than 1 course at a time and track the progress of each course.",
"",
"Interactive learning: The courses within the platform should include various tools for learning such as forums for posting questions and a blog with an overview of the course.",
"",
"",
""
)

$tb.TextFrame.TextRange.Text = [string]::Join("`r", $lines)
$tb.TextFrame.WordWrap = $true
